$d = $word.ActiveDocument

# Locate the "Menu Favorit" list paragraph (the anchor point from the diff:
# the new "Promo" bullet is inserted immediately after it, before "Footer").
$anchor = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text.Trim() -eq "Menu Favorit") {
        $anchor = $p
    }
}

if ($anchor -eq $null) {
    throw "Could not find 'Menu Favorit' paragraph"
}

# Create a new paragraph right after it. InsertParagraphAfter() clones the
# paragraph's pPr (pStyle ListParagraph + numPr ilvl=0/numId=1) and rPr onto
# the freshly created (empty) paragraph, matching the target formatting.
$anchor.Range.InsertParagraphAfter()

# Re-walk the paragraph collection to get a live handle on the paragraph we
# just created (it is the one immediately following the anchor and is empty).
$newPara = $null
$prevWasAnchor = $false
foreach ($p in $d.Paragraphs) {
    if ($prevWasAnchor) {
        $newPara = $p
        $prevWasAnchor = $false
    }
    if ($p.Range.Text.Trim() -eq "Menu Favorit") {
        $prevWasAnchor = $true
    }
}

if ($newPara -eq $null) {
    throw "Could not locate newly inserted paragraph"
}

$target = $newPara.Range

# Inject the two runs exactly as authored ("Promo" + " (Hanya di momen
# tertentu)") via a raw OOXML fragment so they remain distinct <w:r> runs
# (otherwise adjacent runs sharing identical formatting get coalesced into
# one run on save). Paragraph formatting (pStyle/numPr) is preserved because
# it's already set on $target's owning paragraph, but we restate it here too
# so the inserted fragment stays self-consistent if the host reconstructs
# the paragraph mark from the fragment.
$xml = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr><w:rPr><w:rFonts w:ascii="Poppins" w:hAnsi="Poppins" w:cs="Poppins"/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Poppins" w:hAnsi="Poppins" w:cs="Poppins"/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t>Promo</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Poppins" w:hAnsi="Poppins" w:cs="Poppins"/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t xml:space="preserve"> (Hanya di momen tertentu)</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'

$target.InsertXML($xml)
